$d = $word.ActiveDocument
$header = $d.Sections(1).Headers(1)
$header.Range.Text = "Questionnaire 55"
